# Update "想去人数" (want-to-go count) values in F/G columns across the
# three affected sheets: 展览, 演出, 全部类型.
# (本地生活 is untouched by this edit.)

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3019
$ws.Range("F8").Value = 1735
$ws.Range("F9").Value = 1735
$ws.Range("F11").Value = 866
$ws.Range("F12").Value = 140
$ws.Range("F13").Value = 20
$ws.Range("F14").Value = 29
$ws.Range("F15").Value = 2682
$ws.Range("F18").Value = 7184
$ws.Range("F20").Value = 7329
$ws.Range("F23").Value = 5652
$ws.Range("F24").Value = 5652
$ws.Range("F27").Value = 4
$ws.Range("F29").Value = 248
$ws.Range("F30").Value = 202
$ws.Range("F31").Value = 1952
$ws.Range("F36").Value = 499
$ws.Range("F38").Value = 2472
$ws.Range("F39").Value = 1272
$ws.Range("F40").Value = 2867
$ws.Range("F41").Value = 71
$ws.Range("F45").Value = 1121
$ws.Range("F48").Value = 544

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 49
$ws.Range("F12").Value = 363
$ws.Range("G15").Value = 480
$ws.Range("F19").Value = 71

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 49
$ws.Range("F5").Value = 3019
$ws.Range("F7").Value = 1735
$ws.Range("F8").Value = 1735
$ws.Range("F11").Value = 866
$ws.Range("F12").Value = 140
$ws.Range("F13").Value = 20
$ws.Range("F16").Value = 2682
$ws.Range("F20").Value = 7184
$ws.Range("F22").Value = 7329
$ws.Range("F24").Value = 5652
$ws.Range("F25").Value = 5653
$ws.Range("F26").Value = 3144
$ws.Range("F28").Value = 4
$ws.Range("F30").Value = 248
$ws.Range("G31").Value = 480
$ws.Range("F32").Value = 1952
$ws.Range("F37").Value = 499
$ws.Range("F39").Value = 2472
$ws.Range("F40").Value = 1272
$ws.Range("F41").Value = 71
$ws.Range("F42").Value = 2867
$ws.Range("F43").Value = 71
$ws.Range("F48").Value = 1121
$ws.Range("F51").Value = 544
